$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: BASELINE_SIGLA
# Replace existing row 2 and add new rows 3-6
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("BASELINE_SIGLA")

$baseline = @(
    @("MA0304", "Homologation", "BACKUP 16GB RAM 16vCPU (WEBSERVER)", 230, 16, 16, "WEBSERVER", "Linux Server"),
    @("MA0305", "Homologation", "BACKUP 16GB RAM 16vCPU (WEBSERVER)", 230, 16, 16, "WEBSERVER", "Linux Server"),
    @("MA0306", "Production",   "BACKUP 16GB RAM 16vCPU (WEBSERVER)", 230, 16, 16, "WEBSERVER", "Linux Server"),
    @("MA0309", "Production",   "BACKUP 16GB RAM 16vCPU (WEBSERVER)", 230, 16, 16, "WEBSERVER", "Linux Server"),
    @("MA0315", "Development",  "BACKUP 16GB RAM 16vCPU (BACKUP)",      0, 16, 16, "BACKUP",    "Linux Server")
)

$r = 2
foreach ($row in $baseline) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
    $ws1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet 2: Controle
# Overwrite row 2 and add new rows 3-4
# Columns G, H, L, P hold numbers stored as text in the source workbook
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Controle")

$ws2.Cells.Item(2, 1).Value = "Homologation"
$ws2.Cells.Item(2, 2).Value = "4 x BACKUP 32GB RAM 16vCPU (WEBSERVER)"
$ws2.Cells.Item(2, 3).Value = "sa-east-1"
$ws2.Cells.Item(2, 4).Value = "Linux"
$ws2.Cells.Item(2, 5).Value = "c6a.4xlarge"
$ws2.Cells.Item(2, 6).Value = "Shared Instances"
$ws2.Cells.Item(2, 7).Value = "'4"
$ws2.Cells.Item(2, 8).Value = "'40"
$ws2.Cells.Item(2, 9).Value = "Hours/Week"
$ws2.Cells.Item(2, 10).Value = "On-Demand"
$ws2.Cells.Item(2, 11).Value = "General Purpose SSD (gp3)"
$ws2.Cells.Item(2, 12).Value = "'230"
$ws2.Cells.Item(2, 15).Value = "2x Daily"
$ws2.Cells.Item(2, 16).Value = "'10"

$ws2.Cells.Item(3, 1).Value = "Production"
$ws2.Cells.Item(3, 2).Value = "6 x BACKUP 32GB RAM 16vCPU (WEBSERVER)"
$ws2.Cells.Item(3, 3).Value = "sa-east-1"
$ws2.Cells.Item(3, 4).Value = "Linux"
$ws2.Cells.Item(3, 5).Value = "c6a.4xlarge"
$ws2.Cells.Item(3, 6).Value = "Shared Instances"
$ws2.Cells.Item(3, 7).Value = "'6"
$ws2.Cells.Item(3, 9).Value = "Always On"
$ws2.Cells.Item(3, 10).Value = "1 Yr No Upfront EC2 Instance Savings Plan"
$ws2.Cells.Item(3, 11).Value = "General Purpose SSD (gp3)"
$ws2.Cells.Item(3, 12).Value = "'230"
$ws2.Cells.Item(3, 15).Value = "6x Daily"
$ws2.Cells.Item(3, 16).Value = "'20"

$ws2.Cells.Item(4, 1).Value = "Development"
$ws2.Cells.Item(4, 2).Value = "2 x BACKUP 32GB RAM 16vCPU (BACKUP)"
$ws2.Cells.Item(4, 3).Value = "sa-east-1"
$ws2.Cells.Item(4, 4).Value = "Linux"
$ws2.Cells.Item(4, 5).Value = "c6a.4xlarge"
$ws2.Cells.Item(4, 6).Value = "Shared Instances"
$ws2.Cells.Item(4, 7).Value = "'2"
$ws2.Cells.Item(4, 8).Value = "'40"
$ws2.Cells.Item(4, 9).Value = "Hours/Week"
$ws2.Cells.Item(4, 10).Value = "On-Demand"
$ws2.Cells.Item(4, 11).Value = "General Purpose SSD (gp3)"
$ws2.Cells.Item(4, 12).Value = "'0"
$ws2.Cells.Item(4, 15).Value = "2x Daily"
$ws2.Cells.Item(4, 16).Value = "'10"
